# Deployed a4eca05 with MkDocs version: 1.1.2
#
# Repositions a few shapes on slide 3 (title slide) and slide 50 (thank-you
# slide), and bumps a couple of font sizes on slide 50.
#
# NOTE: shape Left/Top/Width/Height are COM "Single" (32-bit float) values
# expressed in points, and the host truncates (rather than rounds) when it
# converts back to EMU on save. The literals below were chosen so that the
# point value, once narrowed to float32, truncates to exactly the target
# EMU from the target OOXML (1 pt = 12700 EMU).

$p = $ppt.ActivePresentation

# ----- Slide 3 -----------------------------------------------------------
$s3 = $p.Slides.Item(3)

# "Straight Connector 7" - nudge the connector's start position
$conn = $s3.Shapes.Item(1)
$conn.Left = 401.6040344238281
$conn.Top = 47.91409683227539

# "TextBox 1" (the rotated "Dan" signature text)
$danBox = $s3.Shapes.Item(3)
$danBox.Left = 406.4718322753906
$danBox.Top = 6.2851972579956055

# ----- Slide 50 -----------------------------------------------------------
$s50 = $p.Slides.Item(50)

# "TextBox 1" - LinkedIn URL: move/resize box and bump font size to 20pt
$linkedin = $s50.Shapes.Item(3)
$linkedin.Left = 509.8086853027344
$linkedin.Top = 216.06039428710938
$linkedin.Width = 391.43267822265625
$linkedin.Height = 31.504724502563477
$linkedin.TextFrame.TextRange.Font.Size = 20

# "Picture 6" - the Twitter/X icon picture: move/resize
$twitterPic = $s50.Shapes.Item(7)
$twitterPic.Left = 460.1380615234375
$twitterPic.Width = 43.255592346191406
$twitterPic.Height = 47.80882263183594

# "TextBox 7" - @dmccreary handle: move/resize box and bump font size to 24pt
$handle = $s50.Shapes.Item(8)
$handle.Left = 509.8086853027344
$handle.Top = 324.189697265625
$handle.Width = 140.18080139160156
$handle.Height = 36.35157775878906
$handle.TextFrame.TextRange.Font.Size = 24
